# Revert "Test: Skill System"
#
# That commit had inserted two new columns ("SkillName" and "Description")
# between "SkillStatName" and "SkillTriggerType" on the skill table sheet.
# Reverting it means removing those two columns again so the sheet goes
# back from A:I to A:G, with everything to the right shifting left by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and D currently hold the "SkillName" / "Description" headers
# (the rest of those columns is empty). Deleting them shifts
# SkillTriggerType/SkillEffectType/SkillAnim/SkillEffectAsset/SkillIcon
# back into C:G.
$ws.Range("C:D").Delete() | Out-Null

# Leave the selection where it naturally ends up after the shift - the
# lone formatted blank cell that was I2 is now G2.
$ws.Range("G2").Select() | Out-Null
